# Update cryptos list values to reflect latest scrape (GitHub Actions run)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.361.20"
$ws.Range("E2").Value = "  +5.72%  "
$ws.Range("D3").Value = "2.509.75"
$ws.Range("E3").Value = "  +3.66%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "324.57"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "105.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +4.33%  "
$ws.Range("E7").Value = "  +2.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("E9").Value = "  +2.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.61%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0817"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.21%  "
$ws.Range("E12").Value = "  +0.82%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "18.36"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.91%  "
$ws.Range("E14").Value = "  +3.74%  "
$ws.Range("D15").Value = "2.902.94"
$ws.Range("E15").Value = "  +3.59%  "
$ws.Range("D16").Value = "2.472.90"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.849"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "
$ws.Range("D18").Value = "47.303.92"
$ws.Range("E18").Value = "  +5.95%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.66"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.87%  "
$ws.Range("E20").Value = "  +2.18%  "
$ws.Range("D21").Value = "0.0₃0941"
$ws.Range("E21").Value = "  +2.57%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "71.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.31%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "252.96"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.51%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.46"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.47%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.90"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.02%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.21"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.47%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "35.25"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +5.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.135"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.58%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.74"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.72%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.75"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.20%  "
$ws.Range("E34").Value = "  +3.24%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.0775"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.19%  "
$ws.Range("E36").Value = "  +0.03%  "
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.77%  "
$ws.Range("B38").Value = "RenderToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.62"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.89%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.98"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +4.62%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "123.23"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.01%  "
$ws.Range("E41").Value = "  +1.99%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.99%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.88"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0296"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("D45").Value = "1.986.22"
$ws.Range("E45").Value = "  +2.27%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.01"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.12"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.46"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +18.41%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "9.02"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.93%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +6.32%  "
